$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.135.98"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "3.310.36"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'599.80"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'143.90"
$ws.Range("E6").Value = "  +5.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.304.90"
$ws.Range("E8").Value = "  +5.93%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "'5.49"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("D12").Value = "'0.474"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'35.02"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "3.856.42"
$ws.Range("E15").Value = "  +5.98%  "
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "3.316.04"
$ws.Range("E17").Value = "  +5.77%  "
$ws.Range("D18").Value = "64.190.48"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'6.91"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").Value = "'484.42"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'14.32"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "'0.745"
$ws.Range("E22").Value = "  +6.48%  "
$ws.Range("D23").Value = "'8.06"
$ws.Range("E23").Value = "  +4.83%  "
$ws.Range("D24").Value = "'13.58"
$ws.Range("D25").Value = "'84.56"
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("B28").Value = "NEARProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'8.30"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'2.16"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").Value = "'28.53"
$ws.Range("E32").Value = "  +4.79%  "
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "'2.58"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").Value = "'6.01"
$ws.Range("E36").Value = "  +2.79%  "
$ws.Range("D37").Value = "'53.43"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").Value = "0.0₃0738"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("E39").Value = "  +2.86%  "
$ws.Range("D40").Value = "'432.09"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Value = "3.013.83"
$ws.Range("E41").Value = "  +4.20%  "
$ws.Range("D42").Value = "'8.47"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").Value = "'2.78"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("E44").Value = "  -5.36%  "
$ws.Range("E45").Value = "  +4.09%  "
$ws.Range("D46").Value = "'2.24"
$ws.Range("E46").Value = "  +5.05%  "
$ws.Range("D47").Value = "'26.41"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "'35.63"
$ws.Range("E50").Value = "  +15.25%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.115"
$ws.Range("E51").Value = "  +1.43%  "
